$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(18.91491841113799, 18.31392817680902, 17.93549525465844, 17.77912295581067, 17.7530336602512, 17.93339480742507, 18.70977401060531, 20.14910027797543, 21.14589348913645, 21.58441245985475, 21.74819509582464, 21.71302464548932, 21.59793323374797, 21.52713658884095, 21.11692354978585, 20.86135177960076, 20.71295533514533, 20.66247484307321, 20.88870352914691, 21.631801064858, 22.10414843704799, 21.85330431492748, 20.87634235589037, 19.76972465227881)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @(8.934646501143355, 8.553639185446668, 8.30971540045676, 8.207902639933087, 8.190854204261546, 8.308351942563723, 8.805402187622111, 9.697246936768098, 10.29791472249717, 10.55867527605942, 10.65558005546374, 10.63479231789674, 10.5666847850996, 10.52472620372483, 10.28061821397793, 10.12763408288384, 10.03846851759531, 10.00807856747031, 10.14404126307671, 10.58673985928057, 10.86533099663491, 10.71763690784685, 10.13662735161579, 9.465311364035953)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colD = @(8.255381932655293, 8.248553358349806, 8.245132061207691, 8.24393306013755, 8.24374579101697, 8.245115099118779, 8.252868063862104, 8.274138136974601, 8.293387824673513, 8.302914214878539, 8.306630672125207, 8.305825444816689, 8.303217794804791, 8.301634681454495, 8.292780599729348, 8.287544870600696, 8.284605959948788, 8.283623403671186, 8.288094726798887, 8.303980779516156, 8.314997874977486, 8.309060335511052, 8.287845914984421, 8.267742109468399)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$colE = @(12.38757432589063, 12.42236024627969, 12.4451525064107, 12.45480160939279, 12.45642566087164, 12.44528117483152, 12.39927133039409, 12.32039421650001, 12.26932443426499, 12.24757791702861, 12.23955609640403, 12.24127426834521, 12.24691368864826, 12.2503957358686, 12.27077546859635, 12.28365787117881, 12.291207328144, 12.29378747697715, 12.28227204806392, 12.24525147418075, 12.22229835810298, 12.23443538571982, 12.28289813263023, 12.34052151164523)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(33.69482902560634, 33.82991092919634, 33.92187426758482, 33.96161145667694, 33.96834614808487, 33.92240103093038, 33.73952826208894, 33.45282663768925, 33.28645215165177, 33.22046178794699, 33.19687360181626, 33.20189132640731, 33.21849305850449, 33.22884474098317, 33.29096044843546, 33.3315543719741, 33.35581499631734, 33.3641856602067, 33.32713863217816, 33.21357865302402, 33.14752931260674, 33.18203145883854, 33.32913211468018, 33.52264338063225)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colI = @(24.77953575933487, 24.93851085218991, 25.04180605976765, 25.0853296821204, 25.09264313726597, 25.04238724358955, 24.83317111780611, 24.46796438094956, 24.2270738735035, 24.12343412451679, 24.08504278387209, 24.09327302402257, 24.12025851283458, 24.13689923136725, 24.23396657297815, 24.29503641491029, 24.33072149140822, 24.34289990296447, 24.28847753712066, 24.11230901947842, 24.00215523971755, 24.06049036950455, 24.29144101399576, 24.56194308863822)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $colI[$i]
}

$colJ = @(9.623628999666686, 9.64899625266407, 9.665372723096034, 9.67224826033255, 9.673402155735683, 9.665464630331805, 9.632209772939254, 9.573323843198493, 9.533878090855264, 9.51675382584034, 9.510386571330375, 9.511752663873196, 9.516227639609728, 9.518983955311858, 9.535013650970969, 9.545056926354379, 9.550910749540924, 9.55290602907375, 9.543979816987532, 9.514910050770844, 9.496594930592568, 9.506307682613469, 9.544466529483033, 9.588580795133092)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $colJ[$i]
}

$colL = @(11.53494986194531, 11.50412389053924, 11.48640380535873, 11.47949057106103, 11.4783613477644, 11.48630931911178, 11.52407314941298, 11.60751024282131, 11.67425736422477, 11.70574246609774, 11.71782070112405, 11.71521261431793, 11.70673306672036, 11.70155917705179, 11.67222181789633, 11.65450750014825, 11.64442461047126, 11.64102908127452, 11.6563822986307, 11.70921954166622, 11.7446553746864, 11.72566186806247, 11.65553438624496, 11.58396128793392)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}

$colO = @(25.83468950717231, 25.95752624055029, 26.03951817855495, 26.07457817576491, 26.08049924515217, 26.03998434463247, 25.8756774466936, 25.60579458203551, 25.4396738810173, 25.37114145514657, 25.34620657256542, 25.35153144160323, 25.36906963955259, 25.37994486464943, 25.4442946603559, 25.48557701027143, 25.50998355900397, 25.51836074916314, 25.48111387941864, 25.36389060886557, 25.29320774303743, 25.33038834964443, 25.48312956556341, 25.67317611036759)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Cells.Item($i + 2, 15).Value = $colO[$i]
}
